# "Actualizar" automatic-update macro: refreshes the availability timestamp
# column (D) by shifting the previously recorded "last checked" timestamps
# down one generation and stamping the most-recent check run at the top.
#
# The sheet is organised in repeating blocks of 14 data rows (one per
# monitored service), each block holding the timestamp of one check run:
#   rows 2-15   -> most recent run
#   rows 16-29  -> previous run
#   rows 30-43  -> run before that
#
# On every refresh:
#   - the oldest block (30-43) is dropped and replaced with a verbatim copy
#     of the previous block (16-29)'s timestamp;
#   - the middle block (16-29) becomes the top block (2-15)'s timestamp,
#     rounded to the nearest millisecond (matching how the live "just
#     checked" value gets archived once it leaves the volatile top row);
#   - the top block (2-15) is stamped with the timestamp of this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockSize  = 14
$firstRow   = 2
$numBlocks  = 3
$col        = 4   # column D
$msPerDay   = 86400000.0

# New timestamp for this refresh run (serial date matching the
# "Actualizar 04-17-2021 13-28-17" automated run).
$newTimestamp = 44303.56099735529

function Round-ToMillisecond([double]$value) {
    $ms = $value * $msPerDay
    $msRounded = [Math]::Floor($ms + 0.5)
    return $msRounded / $msPerDay
}

# Capture the current value of each block's first row (all rows within a
# block share the same timestamp) before overwriting anything.
$blockValues = @()
for ($b = 0; $b -lt $numBlocks; $b++) {
    $row = $firstRow + ($b * $blockSize)
    $blockValues += $ws.Cells.Item($row, $col).Value2
}

# Oldest archived block: verbatim copy of the block one generation newer.
$bottomStart = $firstRow + (($numBlocks - 1) * $blockSize)
$bottomEnd   = $bottomStart + $blockSize - 1
$ws.Range($ws.Cells.Item($bottomStart, $col), $ws.Cells.Item($bottomEnd, $col)).Value2 = $blockValues[$numBlocks - 2]

# Middle block: previous top-block value, archived with millisecond rounding.
$midStart = $firstRow + $blockSize
$midEnd   = $midStart + $blockSize - 1
$roundedTop = Round-ToMillisecond $blockValues[0]
$ws.Range($ws.Cells.Item($midStart, $col), $ws.Cells.Item($midEnd, $col)).Value2 = $roundedTop

# Top block: freshly captured timestamp for this run.
$topStart = $firstRow
$topEnd   = $firstRow + $blockSize - 1
$ws.Range($ws.Cells.Item($topStart, $col), $ws.Cells.Item($topEnd, $col)).Value2 = $newTimestamp
